$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-22 18:48:14'
$ws.Range('E3').Value = '2026-02-22 18:48:17'
$ws.Range('K3').Value = '15.6 MJ/m2'
$ws.Range('O3').Value = '4.7 °C'
$ws.Range('E4').Value = '2026-02-22 18:48:19'
$ws.Range('H4').Value = "'62%"
$ws.Range('E5').Value = '2026-02-22 18:48:22'
$ws.Range('E6').Value = '2026-02-22 18:48:24'
$ws.Range('E7').Value = '2026-02-22 18:48:27'
$ws.Range('E8').Value = '2026-02-22 18:48:29'
$ws.Range('H8').Value = "'49%"
$ws.Range('K8').Value = '16.0 MJ/m2'
$ws.Range('E9').Value = '2026-02-22 18:48:31'
$ws.Range('H9').Value = "'74%"
$ws.Range('E10').Value = '2026-02-22 18:48:32'
$ws.Range('H10').Value = "'79%"
$ws.Range('E11').Value = '2026-02-22 18:48:34'
$ws.Range('E12').Value = '2026-02-22 18:48:35'
$ws.Range('E13').Value = '2026-02-22 18:48:36'
$ws.Range('H13').Value = "'61%"
$ws.Range('J13').Value = '1030.4 hPa'
$ws.Range('L13').Value = '21.2 km/h - 93º 18:05 TU'
$ws.Range('O13').Value = '6.4 °C'
$ws.Range('E14').Value = '2026-02-22 18:48:37'
$ws.Range('E15').Value = '2026-02-22 18:48:38'
$ws.Range('H15').Value = "'71%"
$ws.Range('E16').Value = '2026-02-22 18:48:39'
$ws.Range('E17').Value = '2026-02-22 18:48:40'
$ws.Range('H17').Value = "'28%"
$ws.Range('E18').Value = '2026-02-22 18:48:41'
$ws.Range('J18').Value = '1027.6 hPa'
$ws.Range('K18').Value = '15.2 MJ/m2'
$ws.Range('O18').Value = '10.2 °C'
$ws.Range('E19').Value = '2026-02-22 18:48:42'
$ws.Range('E20').Value = '2026-02-22 18:48:43'
$ws.Range('O20').Value = '4.0 °C'
$ws.Range('E21').Value = '2026-02-22 18:48:46'
$ws.Range('J21').Value = '1029.2 hPa'
$ws.Range('E22').Value = '2026-02-22 18:48:48'
$ws.Range('H22').Value = "'22%"
$ws.Range('E23').Value = '2026-02-22 18:48:51'
$ws.Range('L23').Value = '23.8 km/h - 329º 18:08 TU'
$ws.Range('O23').Value = '5.9 °C'
$ws.Range('E24').Value = '2026-02-22 18:48:53'
$ws.Range('O24').Value = '7.8 °C'
$ws.Range('E25').Value = '2026-02-22 18:48:56'
$ws.Range('O25').Value = '7.3 °C'
$ws.Range('E26').Value = '2026-02-22 18:48:58'
$ws.Range('H26').Value = "'34%"
$ws.Range('J26').Value = '1026.0 hPa'
$ws.Range('O26').Value = '11.6 °C'
$ws.Range('E27').Value = '2026-02-22 18:49:01'
$ws.Range('E28').Value = '2026-02-22 18:49:03'
$ws.Range('E29').Value = '2026-02-22 18:49:05'
$ws.Range('E30').Value = '2026-02-22 18:49:08'
$ws.Range('H30').Value = "'70%"
$ws.Range('J30').Value = '1027.2 hPa'
$ws.Range('E31').Value = '2026-02-22 18:49:10'
$ws.Range('J31').Value = '1026.6 hPa'
$ws.Range('O31').Value = '14.4 °C'
$ws.Range('E32').Value = '2026-02-22 18:49:13'
$ws.Range('O32').Value = '6.4 °C'
$ws.Range('E33').Value = '2026-02-22 18:49:15'
$ws.Range('J33').Value = '1028.7 hPa'
$ws.Range('O33').Value = '8.3 °C'
$ws.Range('E34').Value = '2026-02-22 18:49:18'
$ws.Range('H34').Value = "'44%"
$ws.Range('O34').Value = '4.6 °C'
$ws.Range('E35').Value = '2026-02-22 18:49:20'
$ws.Range('H35').Value = "'41%"
$ws.Range('E36').Value = '2026-02-22 18:49:23'
$ws.Range('O36').Value = '11.8 °C'
$ws.Range('E37').Value = '2026-02-22 18:49:25'
$ws.Range('J37').Value = '1029.8 hPa'
$ws.Range('E38').Value = '2026-02-22 18:49:27'
$ws.Range('E39').Value = '2026-02-22 18:49:30'
$ws.Range('E40').Value = '2026-02-22 18:49:32'
$ws.Range('E41').Value = '2026-02-22 18:49:35'
$ws.Range('J41').Value = '1027.7 hPa'
$ws.Range('E42').Value = '2026-02-22 18:49:37'
$ws.Range('E43').Value = '2026-02-22 18:49:39'
$ws.Range('O43').Value = '9.2 °C'
$ws.Range('E44').Value = '2026-02-22 18:49:42'
$ws.Range('E45').Value = '2026-02-22 18:49:44'
$ws.Range('H45').Value = "'53%"
$ws.Range('O45').Value = '9.0 °C'
$ws.Range('E46').Value = '2026-02-22 18:49:47'
$ws.Range('O46').Value = '9.0 °C'
